$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the subject label in B3
$ws.Range("B3").Value = "Test"

# Replace the first active student row (row 10) with the new student
$ws.Range("A10").Value = "KONDI"
$ws.Range("B10").Value = "Abdoul Malik"

# Remove the remaining inactive students (rows 11-19)
$ws.Range("A11:P19").EntireRow.Delete()
